# Append new invoice line-item rows (12-15) to Sheet1, matching the
# pattern of the existing data (two new invoices, each with two line items).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("SP22092022132249", 100, "PANADOL STRIP 10", 1, 15),
    @("SP22092022132249", 101, "PANADOL STRIP 20", 1, 29),
    @("SP22092022154655", 100, "PANADOL STRIP 10", 1, 15),
    @("SP22092022154655", 101, "PANADOL STRIP 20", 1, 29)
)

$startRow = 12
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
}
